$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top to keep row indices stable while deleting.
$ws.Rows("10").Delete() | Out-Null
$ws.Rows("9").Delete() | Out-Null
$ws.Rows("8").Delete() | Out-Null
$ws.Rows("7").Delete() | Out-Null
$ws.Rows("6").Delete() | Out-Null
$ws.Rows("3").Delete() | Out-Null

# Update the auth_key value in row 3 (previously "Switch drugi" row) from "Password" to "ok".
$ws.Range("K3").Value = "ok"

# Update the active cell selection to J7.
$ws.Range("J7").Select() | Out-Null
